# Hortaliza, Terminal Hortofrutícola Agro Chillán - Pimiento
# Insert two new weekly price records (Zafiro rojo / Zafiro verde,
# fecha 44504) right above the current row 146, pushing every
# following record down by two rows (dimension grows from R173 to R175).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 146:173 down to make room for the two new records.
$ws.Rows("146:147").Insert()

# New row 146: Pimiento, Zafiro rojo, Región de Arica y Parinacota
$ws.Range("A146").Value = 7
$ws.Range("B146").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C146").Value = "Ñuble"
$ws.Range("D146").Value = 44504
$ws.Range("E146").Value = 16
$ws.Range("F146").Value = 100112002
$ws.Range("G146").Value = "Pimiento"
$ws.Range("H146").Value = "Zafiro rojo"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 100
$ws.Range("K146").Value = 37000
$ws.Range("L146").Value = 38000
$ws.Range("M146").Value = 37500
$ws.Range("N146").Value = "$/caja 15 kilos"
$ws.Range("O146").Value = "Región de Arica y Parinacota"
$ws.Range("P146").Value = 2500
$ws.Range("Q146").Value = 15
$ws.Range("R146").Value = "Hortaliza"

# New row 147: Pimiento, Zafiro verde, Región de Arica y Parinacota
$ws.Range("A147").Value = 7
$ws.Range("B147").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C147").Value = "Ñuble"
$ws.Range("D147").Value = 44504
$ws.Range("E147").Value = 16
$ws.Range("F147").Value = 100112002
$ws.Range("G147").Value = "Pimiento"
$ws.Range("H147").Value = "Zafiro verde"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 100
$ws.Range("K147").Value = 37000
$ws.Range("L147").Value = 38000
$ws.Range("M147").Value = 37500
$ws.Range("N147").Value = "$/caja 15 kilos"
$ws.Range("O147").Value = "Región de Arica y Parinacota"
$ws.Range("P147").Value = 2500
$ws.Range("Q147").Value = 15
$ws.Range("R147").Value = "Hortaliza"
